$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N, pushing the existing
# "Late"/"Waived"/"Outstanding" columns one slot to the right.
$mWidth = $ws.Columns("M").ColumnWidth
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $mWidth

# Make "Repayment schedule" the active sheet / tab, and update the
# selected cell on that sheet.
$ws.Activate()
$ws.Range("S9").Select() | Out-Null
